# feat: add 2022-Q4 data
#
# 1) Insert a new row at the top of the "总计" (totals) sheet for the
#    2022-Q4 quarter, pushing the existing quarters down by one row.
# 2) Insert a new worksheet named "2022-Q4" right before the existing
#    "2021-Q4" sheet (a duplicate of it, carrying over its layout/styling),
#    then update its figures to the new quarter's numbers.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" sheet: add the 2022-Q4 summary row
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")

$totals.Rows.Item(2).Insert()
# New row inherited the row-1 (header) look for B:D - strip it back to the
# plain/no-style look the other data rows use.
$totals.Range("B2:D2").Style = "Normal"
# Column A carries its own bold/border style on every data row - copy it
# (formats only) from the row below rather than guessing the style index.
$totals.Cells.Item(3, 1).Copy()
$totals.Cells.Item(2, 1).PasteSpecial(-4122)

$totals.Cells.Item(2, 1).Value = 4
$totals.Cells.Item(2, 2).Value = "2022-Q4"
$totals.Cells.Item(2, 3).Value = 2
$totals.Cells.Item(2, 4).Value = 0.13

# ---------------------------------------------------------------------
# 2) New "2022-Q4" sheet, cloned from "2021-Q4"
# ---------------------------------------------------------------------
$q4_2021 = $wb.Worksheets.Item("2021-Q4")
$q4_2021.Copy($q4_2021)

$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q4"

function Set-TextValue($cell, $text) {
    $originalStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $originalStyle
}

Set-TextValue $newSheet.Cells.Item(2, 4) "1.79"
Set-TextValue $newSheet.Cells.Item(2, 5) "88.58"
Set-TextValue $newSheet.Cells.Item(2, 7) "0.0626"

Set-TextValue $newSheet.Cells.Item(3, 4) "1.79"
Set-TextValue $newSheet.Cells.Item(3, 5) "88.58"
Set-TextValue $newSheet.Cells.Item(3, 7) "0.0626"

# Creating/renaming the new sheet makes it the active tab - restore the
# original active sheet ("2020-Q4") so the view state isn't disturbed.
$wb.Worksheets.Item("2020-Q4").Activate()

Write-Host "done"
